# The "log" table (row 6) is renamed to "registro" and gets two new
# trailing columns: "data" and "id venda". Columns C6/D6/E6
# (id usuario / tipo / id produto) are untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "data" column first so it is appended to the shared-string
# table ahead of "registro" (matches original authoring order).
$ws.Range("F6").Value = "data"

# Rename the table header itself from "log" to "registro".
$ws.Range("A6").Value = "registro"

# Add the final new column "id venda".
$ws.Range("G6").Value = "id venda"

# Move the active selection to the newly added last cell, G6.
$ws.Range("G6").Select()
